# Arbeitszeit_Pichler.xlsx - Tagebuch: add two new journal entries (rows 45 & 46)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Reuse the date-formatted style already used by the column (E6:E44)
# by copying formats from the last existing entry row (44) down onto
# the two new rows before writing values.
$ws.Cells.Item(44, 5).Copy()
$ws.Range("E45:E46").PasteSpecial(-4122)

# Row 45: 25.08.2019, 1h, Einheit "Stunden", Taetigkeit "Statusbericht"
$ws.Cells.Item(45, 5).Value = 43702
$ws.Cells.Item(45, 6).Value = 1
$ws.Cells.Item(45, 7).Value = "Stunden"
$ws.Cells.Item(45, 8).Value = "Statusbericht"

# Row 46: 27.08.2019, 1h, Einheit "Stunden", Taetigkeit "Programmieren"
$ws.Cells.Item(46, 5).Value = 43704
$ws.Cells.Item(46, 6).Value = 1
$ws.Cells.Item(46, 7).Value = "Stunden"
$ws.Cells.Item(46, 8).Value = "Programmieren"

# Keep selection consistent with the new last used cell
$ws.Range("H46").Select()
